$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.184.24'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '1.833.39'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('D4').Value = '''0.9990'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''241.85'
$ws.Range('D6').Value = '''0.6642'
$ws.Range('E6').Value = '  -2.29%  '
$ws.Range('D7').Value = '''0.9999'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '''0.07416'
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').Value = '''0.2935'
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').Value = '''22.90'
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').Value = '''0.07756'
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').Value = '1.837.41'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = '''4.989'
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').Value = '''0.6682'
$ws.Range('E14').Value = '  -1.06%  '
$ws.Range('D15').Value = '''82.89'
$ws.Range('E15').Value = '  -4.22%  '
$ws.Range('D16').Value = '''6.097'
$ws.Range('E16').Value = '  -0.82%  '
$ws.Range('D17').Value = '''0.000008360'
$ws.Range('E17').Value = '  +1.56%  '
$ws.Range('D18').Value = '29.190.01'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').Value = '2.089.12'
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '''228.21'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').Value = '''7.149'
$ws.Range('E23').Value = '  -2.46%  '
$ws.Range('D24').Value = '''0.9998'
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '''159.41'
$ws.Range('E25').Value = '  -0.93%  '
$ws.Range('D26').Value = '''0.1413'
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('D27').Value = '''8.611'
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('D29').Value = '''1.513'
$ws.Range('E29').Value = '  +0.86%  '
$ws.Range('E30').Value = '  -3.15%  '
$ws.Range('D31').Value = '''4.039'
$ws.Range('E31').Value = '  -2.13%  '
$ws.Range('D32').Value = '''1.191'
$ws.Range('E32').Value = '  -0.32%  '
$ws.Range('D33').Value = '''0.05304'
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('D34').Value = '''1.866'
$ws.Range('E34').Value = '  +1.07%  '
$ws.Range('D35').Value = '''0.7473'
$ws.Range('E35').Value = '  -0.82%  '
$ws.Range('E36').Value = '  +1.07%  '
$ws.Range('D37').Value = '''2.646'
$ws.Range('E37').Value = '  -1.40%  '
$ws.Range('D38').Value = '1.273.03'
$ws.Range('E38').Value = '  -2.52%  '
$ws.Range('E39').Value = '  -0.82%  '
$ws.Range('D40').Value = '''2.734'
$ws.Range('E40').Value = '  +0.74%  '
$ws.Range('D41').Value = '''0.9334'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = '''5.887'
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('D43').Value = '''0.08433'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('D45').Value = '''101.93'
$ws.Range('E45').Value = '  -2.88%  '
$ws.Range('D46').Value = '1.990.50'
$ws.Range('E46').Value = '  +0.34%  '
$ws.Range('D47').Value = '''0.5148'
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('D48').Value = '''1.760'
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').Value = '''62.99'
$ws.Range('E50').Value = '  -1.17%  '
$ws.Range('D51').Value = '''0.05877'
$ws.Range('E51').Value = '  -0.87%  '
